$d = $word.ActiveDocument

# 1. Update table cell margin (left) from 118 dxa (5.9pt) to 123 dxa (6.15pt)
$t = $d.Tables(1)
$t.LeftPadding = 6.15

# 2. Merge the date + "г." runs into a single run with updated placeholder text
$found = $d.Content.Find.Execute(
    "{issue.humanized_created_at_with_month_as_word} г.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{issue.humanized_created_at_with_quotes_and_month_as_word} г.", 2
)

# 3. Add the two new character styles (ListLabel 9 / ListLabel 10)
$s9 = $d.Styles.Add("ListLabel 9", 2)
$s9.Font.Name = "Times New Roman"
$s9.Font.NameAscii = "Times New Roman"
$s9.Font.Bold = $true
$s9.Font.Size = 10.5
$s9.QuickStyle = $true

$s10 = $d.Styles.Add("ListLabel 10", 2)
$s10.Font.Name = "Times New Roman"
$s10.Font.NameAscii = "Times New Roman"
$s10.Font.Bold = $false
$s10.Font.Size = 10.5
$s10.QuickStyle = $true
